$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2021 declaration) updates ---
# id: z0bug.li_partner_6 -> z0bug.li_partner_6_2021
$ws.Range("A2").Value = "z0bug.li_partner_6_2021"

# date / partner_document_date / date_start: <#-01-15 -> <###-01-25
$ws.Range("E2:G2").Value = "<###-01-25"

# date_end: <#-12-31 -> <###-12-31
$ws.Range("H2").Value = "<###-12-31"

# --- Row 3 (2022 declaration) updates ---
# id: z0bug.li_partner_6 -> z0bug.li_partner_6_2022
$ws.Range("A3").Value = "z0bug.li_partner_6_2022"

# date / partner_document_date / date_start: ####-01-20 -> ####-01-06
$ws.Range("E3:G3").Value = "####-01-06"

# --- Cosmetic / view changes ---
# Widen column A and move the active selection to A4
$ws.Range("A1").EntireColumn.ColumnWidth = 21.1
$null = $ws.Range("A4").Select()
